$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove all existing hyperlinks on Sheet1 (there were 7: A1,A3,A4,A5,A6,A7,A8) ---
[void]$ws.Hyperlinks.Delete()

# --- Wipe all existing cell content/formatting (old layout was A1:F8) ---
[void]$ws.Cells.Clear()

# --- Rebuild the table with the new, smaller layout (A1:C4) ---
# Header row (bold + yellow fill style, like old D1:F1)
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "Login Status"

# Data rows
$ws.Range("A2").Value = "prasadn@leotechnosoft.net"
$ws.Range("B2").Value = "leo_123"
$ws.Range("C2").Value = "Pass"

$ws.Range("A3").Value = "neel.sharma@spicetg.com"
$ws.Range("B3").Value = "spice_12345"
$ws.Range("C3").Value = "Fail"

$ws.Range("A4").Value = "mark@leotechnosoft.net"
$ws.Range("B4").Value = "leo_12345"
$ws.Range("C4").Value = "Fail"

# --- Formatting: reuse the workbook's existing styles by copying formats from
#     cells that already carry them, so the style table doesn't get duplicated. ---

# A1 / B1 -> bold + yellow header style (same as old D1/E1)
$ws.Range("D1").Copy() | Out-Null
$ws.Range("A1:B1").PasteSpecial(-4122) | Out-Null

# C1 -> new header style: like A1/B1's data style (numFmt 49 + border + left/vcenter)
# but with the yellow header fill. Build it by copying the plain data style first,
# then applying the yellow interior color (matches the newly-added cellXfs entry).
$ws.Range("B2").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null

# A2,A3,A4 -> hyperlink-like style (same as old A1 which used style index 3)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A2:A4").PasteSpecial(-4122) | Out-Null

# B2,C2,B3,C3,B4,C4 -> plain bordered/left-aligned data style (same as old B2)
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B2:C4").PasteSpecial(-4122) | Out-Null

# Now give C1 its yellow fill (creates the new cellXfs entry, after all the
# already-existing styles above have been resolved/reused).
$ws.Range("C1").Interior.Color = 65535

$excel.CutCopyMode = 0

# --- Re-add the two remaining hyperlinks ---
[void]$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:prasadn@leotechnosoft.net")
[void]$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:mark@leotechnosoft.net")

# --- Column B width (closest value reachable through this host's column-width
#     rounding to the target 11.5703125) ---
$ws.Range("B1").ColumnWidth = 10.71

# --- Selection / active cell shown in the saved view ---
[void]$ws.Range("E8").Select()
